$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's columns D,E,F,G were mislabelled:
#   D = codeforiati:category-code   E = codeforiati:group-code
#   F = codeforiati:group-name      G = codeforiati:category-name
# The fix swaps the (D,E) pair of cells with the (F,G) pair of cells in
# every row (header + data), i.e.:
#   new D = old F     new E = old G
#   new F = old E     new G = old D
# Range.Copy (rather than .Value) is used so the shared-string cell type
# is preserved instead of being coerced to a number for numeric-looking
# codes like "111".

for ($r = 1; $r -le 235; $r++) {
    # Stash the current D and E values in scratch cells, in the order
    # they are needed later (AA <- D, Z <- E) so a straight two-cell
    # copy back out lands them correctly in F,G.
    $ws.Range("D" + $r).Copy($ws.Range("AA" + $r))
    $ws.Range("E" + $r).Copy($ws.Range("Z" + $r))

    # F,G -> D,E
    $ws.Range("F" + $r + ":G" + $r).Copy($ws.Range("D" + $r))

    # stashed old E,D -> F,G  (Z holds old E, AA holds old D)
    $ws.Range("Z" + $r + ":AA" + $r).Copy($ws.Range("F" + $r))

    $ws.Range("Z" + $r + ":AA" + $r).Clear()
}
